$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pgf"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 22.51188433333333
$ws.Range("H2").Value = 67.535653
$ws.Range("I2").Value = 0.7173237801266834
$ws.Range("J2").Value = 0.7173237801266834
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 52.91030366666666
$ws.Range("N2").Value = 158.730911
$ws.Range("O2").Value = 0.4161415425564564
$ws.Range("P2").Value = 0.4161415425564564
$ws.Range("Q2").Value = 1191.110636185542
$ws.Range("R2").Value = 10719.99572566988
$ws.Range("S2").Value = 0.2985082243743464
$ws.Range("T2").Value = 0.2985082243743464

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pgf"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 22.51188433333333
$ws.Range("H3").Value = 67.535653
$ws.Range("I3").Value = 0.7173237801266834
$ws.Range("J3").Value = 0.7173237801266834
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.26005833333333
$ws.Range("N3").Value = 141.780175
$ws.Range("O3").Value = 0.3717021489810786
$ws.Range("P3").Value = 0.3717021489810786
$ws.Range("Q3").Value = 1063.912966786586
$ws.Range("R3").Value = 9575.216701079273
$ws.Range("S3").Value = 0.2666307905883189
$ws.Range("T3").Value = 0.2666307905883189

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pgf"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 22.51188433333333
$ws.Range("H4").Value = 67.535653
$ws.Range("I4").Value = 0.7173237801266834
$ws.Range("J4").Value = 0.7173237801266834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 26.97460733333333
$ws.Range("N4").Value = 80.923822
$ws.Range("O4").Value = 0.2121563084624651
$ws.Range("P4").Value = 0.2121563084624651
$ws.Range("Q4").Value = 607.249240225085
$ws.Range("R4").Value = 5465.243162025766
$ws.Range("S4").Value = 0.1521847651640182
$ws.Range("T4").Value = 0.1521847651640182

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pgf"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.074492666666667
$ws.Range("H5").Value = 21.223478
$ws.Range("I5").Value = 0.2254232363222357
$ws.Range("J5").Value = 0.2254232363222357
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 52.91030366666666
$ws.Range("N5").Value = 158.730911
$ws.Range("O5").Value = 0.4161415425564564
$ws.Range("P5").Value = 0.4161415425564564
$ws.Range("Q5").Value = 374.3135552809398
$ws.Range("R5").Value = 3368.821997528458
$ws.Range("S5").Value = 0.09380797329120376
$ws.Range("T5").Value = 0.09380797329120374

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pgf"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.074492666666667
$ws.Range("H6").Value = 21.223478
$ws.Range("I6").Value = 0.2254232363222357
$ws.Range("J6").Value = 0.2254232363222357
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 47.26005833333333
$ws.Range("N6").Value = 141.780175
$ws.Range("O6").Value = 0.3717021489810786
$ws.Range("P6").Value = 0.3717021489810786
$ws.Range("Q6").Value = 334.3409361054055
$ws.Range("R6").Value = 3009.06842494865
$ws.Range("S6").Value = 0.08379030137124453
$ws.Range("T6").Value = 0.08379030137124452

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pgf"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.074492666666667
$ws.Range("H7").Value = 21.223478
$ws.Range("I7").Value = 0.2254232363222357
$ws.Range("J7").Value = 0.2254232363222357
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 26.97460733333333
$ws.Range("N7").Value = 80.923822
$ws.Range("O7").Value = 0.2121563084624651
$ws.Range("P7").Value = 0.2121563084624651
$ws.Range("Q7").Value = 190.8316617658796
$ws.Range("R7").Value = 1717.484955892916
$ws.Range("S7").Value = 0.04782496165978741
$ws.Range("T7").Value = 0.0478249616597874

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Pgf"
$ws.Range("C8").Value = "Nrp1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.796779333333333
$ws.Range("H8").Value = 5.390338
$ws.Range("I8").Value = 0.05725298355108089
$ws.Range("J8").Value = 0.05725298355108089
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 52.91030366666666
$ws.Range("N8").Value = 158.730911
$ws.Range("O8").Value = 0.4161415425564564
$ws.Range("P8").Value = 0.4161415425564564
$ws.Range("Q8").Value = 95.06814014865755
$ws.Range("R8").Value = 855.613261337918
$ws.Range("S8").Value = 0.02382534489090622
$ws.Range("T8").Value = 0.02382534489090622

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Pgf"
$ws.Range("C9").Value = "Nrp1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.796779333333333
$ws.Range("H9").Value = 5.390338
$ws.Range("I9").Value = 0.05725298355108089
$ws.Range("J9").Value = 0.05725298355108089
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 47.26005833333333
$ws.Range("N9").Value = 141.780175
$ws.Range("O9").Value = 0.3717021489810786
$ws.Range("P9").Value = 0.3717021489810786
$ws.Range("Q9").Value = 84.9158961054611
$ws.Range("R9").Value = 764.2430649491499
$ws.Range("S9").Value = 0.02128105702151511
$ws.Range("T9").Value = 0.02128105702151511

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pgf"
$ws.Range("C10").Value = "Nrp1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.796779333333333
$ws.Range("H10").Value = 5.390338
$ws.Range("I10").Value = 0.05725298355108089
$ws.Range("J10").Value = 0.05725298355108089
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 26.97460733333333
$ws.Range("N10").Value = 80.923822
$ws.Range("O10").Value = 0.2121563084624651
$ws.Range("P10").Value = 0.2121563084624651
$ws.Range("Q10").Value = 48.46741698131511
$ws.Range("R10").Value = 436.206752831836
$ws.Range("S10").Value = 0.01214658163865956
$ws.Range("T10").Value = 0.01214658163865956

Write-Output "edit complete"
